$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 6059965
$ws.Range("I19").Value = 334453.44
$ws.Range("J19").Value = 12501166
$ws.Range("K19").Value = 334453.44
$ws.Range("L19").Value = 12501166
$ws.Range("M19").Value = -334278.44
$ws.Range("N19").Value = -12501516
$ws.Range("H62").Value = 2904.1052
$ws.Range("I62").Value = 3042.4546
$ws.Range("K62").Value = 3042.4546
$ws.Range("M62").Value = -2418.4546
$ws.Range("H65").Value = 2904.1052
$ws.Range("I65").Value = 3042.4546
$ws.Range("K65").Value = 15212.273
$ws.Range("M65").Value = -12092.273
$ws.Range("H111").Value = 125004060
$ws.Range("I111").Value = 166671870
$ws.Range("J111").Value = 650
$ws.Range("K111").Value = 500015610
$ws.Range("L111").Value = 1950
$ws.Range("M111").Value = -500012543
$ws.Range("N111").Value = -8084
$ws.Range("H116").Value = 2567592.5
$ws.Range("I116").Value = 9617748
$ws.Range("J116").Value = 3899.7727
$ws.Range("K116").Value = 9617748
$ws.Range("L116").Value = 3899.7727
$ws.Range("M116").Value = -9614306
$ws.Range("N116").Value = -10783.7727
$ws.Range("H123").Value = 46488
$ws.Range("J123").Value = 46488
$ws.Range("L123").Value = 46488
$ws.Range("N123").Value = -56288
$ws.Range("H129").Value = 1100.41
$ws.Range("J129").Value = 1115.9485
$ws.Range("L129").Value = 3347.8455
$ws.Range("N129").Value = -13347.8455
$ws.Range("H131").Value = 4127.7144
$ws.Range("I131").Value = 2644.75
$ws.Range("J131").Value = 6105
$ws.Range("K131").Value = 7934.25
$ws.Range("L131").Value = 18315
$ws.Range("M131").Value = -2894.25
$ws.Range("N131").Value = -28395
$ws.Range("H139").Value = 10490.667
$ws.Range("J139").Value = 10490.667
$ws.Range("L139").Value = 10490.667
$ws.Range("N139").Value = -20770.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1314.2894
$ws.Range("I74").Value = 1249.3226
$ws.Range("J74").Value = 1602
$ws.Range("K74").Value = 1249.3226
$ws.Range("L74").Value = 1602
$ws.Range("M74").Value = -375.3226
$ws.Range("N74").Value = -3350
$ws.Range("H77").Value = 1314.2894
$ws.Range("I77").Value = 1249.3226
$ws.Range("J77").Value = 1602
$ws.Range("K77").Value = 6246.612999999999
$ws.Range("L77").Value = 8010
$ws.Range("M77").Value = -1878.612999999999
$ws.Range("N77").Value = -16746
$ws.Range("H113").Value = 79632.664
$ws.Range("J113").Value = 79632.664
$ws.Range("L113").Value = 79632.664
$ws.Range("N113").Value = -88310.664
$ws.Range("H122").Value = 1996.5
$ws.Range("I122").Value = 1995.8
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5987.4
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3537.4
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -27
$ws.Range("N22").Value = -546
$ws.Range("H134").Value = 1887.6774
$ws.Range("I134").Value = 1791.9231
$ws.Range("J134").Value = 2385.6
$ws.Range("K134").Value = 5375.7693
$ws.Range("L134").Value = 7156.799999999999
$ws.Range("M134").Value = -2840.7693
$ws.Range("N134").Value = -12226.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2604.5925
$ws.Range("I31").Value = 2031.5555
$ws.Range("J31").Value = 3750.6667
$ws.Range("K31").Value = 2031.5555
$ws.Range("L31").Value = 3750.6667
$ws.Range("M31").Value = -1736.5555
$ws.Range("N31").Value = -4340.6667
$ws.Range("H34").Value = 2604.5925
$ws.Range("I34").Value = 2031.5555
$ws.Range("J34").Value = 3750.6667
$ws.Range("K34").Value = 2031.5555
$ws.Range("L34").Value = 3750.6667
$ws.Range("M34").Value = -1829.5555
$ws.Range("N34").Value = -4154.6667
$ws.Range("H137").Value = 30884.137
$ws.Range("J137").Value = 30884.137
$ws.Range("L137").Value = 30884.137
$ws.Range("N137").Value = -41084.137
$ws.Range("H138").Value = 39352.47
$ws.Range("J138").Value = 39352.47
$ws.Range("L138").Value = 39352.47
$ws.Range("N138").Value = -49632.47

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1144.4445
$ws.Range("I86").Value = 850
$ws.Range("J86").Value = 1733.3334
$ws.Range("K86").Value = 2550
$ws.Range("L86").Value = 5200.0002
$ws.Range("M86").Value = -1364
$ws.Range("N86").Value = -7572.0002
$ws.Range("H89").Value = 1144.4445
$ws.Range("I89").Value = 850
$ws.Range("J89").Value = 1733.3334
$ws.Range("K89").Value = 7650
$ws.Range("L89").Value = 15600.0006
$ws.Range("M89").Value = -1722
$ws.Range("N89").Value = -27456.0006
$ws.Range("H97").Value = 717.4
$ws.Range("I97").Value = 489.5
$ws.Range("J97").Value = 869.3333
$ws.Range("K97").Value = 1468.5
$ws.Range("L97").Value = 2607.9999
$ws.Range("M97").Value = -972.5
$ws.Range("N97").Value = -3599.9999
$ws.Range("H101").Value = 4800
$ws.Range("H122").Value = 8714.370000000001
$ws.Range("I122").Value = 10878.048
$ws.Range("J122").Value = 1141.5
$ws.Range("K122").Value = 97902.432
$ws.Range("L122").Value = 10273.5
$ws.Range("M122").Value = -95452.432
$ws.Range("N122").Value = -15173.5
$ws.Range("H123").Value = 4001.111
$ws.Range("I123").Value = 2010
$ws.Range("J123").Value = 4996.6665
$ws.Range("K123").Value = 6030
$ws.Range("L123").Value = 14989.9995
$ws.Range("M123").Value = -3580
$ws.Range("N123").Value = -19889.9995
$ws.Range("H131").Value = 770.6129
$ws.Range("I131").Value = 476.66666
$ws.Range("J131").Value = 890.86365
$ws.Range("K131").Value = 1429.99998
$ws.Range("L131").Value = 2672.59095
$ws.Range("M131").Value = 3610.00002
$ws.Range("N131").Value = -12752.59095
$ws.Range("H132").Value = 843204.2
$ws.Range("I132").Value = 904
$ws.Range("J132").Value = 1011664.2
$ws.Range("K132").Value = 8136
$ws.Range("L132").Value = 9104977.799999999
$ws.Range("M132").Value = -5606
$ws.Range("N132").Value = -9110037.799999999
$ws.Range("H133").Value = 4500
$ws.Range("I133").Value = 3000
$ws.Range("K133").Value = 9000
$ws.Range("M133").Value = -3940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 49750
$ws.Range("J141").Value = 49750
$ws.Range("L141").Value = 49750
$ws.Range("N141").Value = -60110

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 391.81818
$ws.Range("I55").Value = 370
$ws.Range("J55").Value = 450
$ws.Range("K55").Value = 370
$ws.Range("L55").Value = 450
$ws.Range("M55").Value = -197
$ws.Range("N55").Value = -796
$ws.Range("H87").Value = 33695
$ws.Range("J87").Value = 33695
$ws.Range("L87").Value = 33695
$ws.Range("N87").Value = -35941
$ws.Range("H88").Value = 25990
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 25990
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 25990
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -26846
$ws.Range("H90").Value = 33695
$ws.Range("J90").Value = 33695
$ws.Range("L90").Value = 101085
$ws.Range("N90").Value = -112317
$ws.Range("H91").Value = 25990
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 25990
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 25990
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -28954
$ws.Range("H122").Value = 2980
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2980
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 8940
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -13840
$ws.Range("H128").Value = 56421.75
$ws.Range("J128").Value = 56421.75
$ws.Range("L128").Value = 56421.75
$ws.Range("N128").Value = -66381.75
$ws.Range("H140").Value = 36995
$ws.Range("J140").Value = 36995
$ws.Range("L140").Value = 36995
$ws.Range("N140").Value = -47355

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 60000
$ws.Range("J46").Value = 60000
$ws.Range("L46").Value = 60000
$ws.Range("N46").Value = -60462
$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 180000
$ws.Range("N134").Value = -185070
